# Generate Report for Handback
#
# The handback report was regenerated, which refreshed the
# "Latest Handback DateTime" value (column K, row 2 - the
# a99ca549-0f60-4d18-ad94-4b495b53a99d.md entry) on both the
# "zh-cn" and "de-de" locale-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-09-07 14:45:03"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-09-07 14:45:52"
